# Update cryptocurrency price/volume figures per the scraper refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.655.14"
$ws.Range("D3").Value = "1.636.58"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'213.23"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").Value = "'0.0623"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "'19.14"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.865.49"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "1.637.66"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "'0.528"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "26.664.50"
$ws.Range("D17").Value = "'63.30"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").Value = "'218.55"
$ws.Range("E19").Value = "  +7.79%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").Value = "'6.23"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "'149.09"
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'6.86"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "'15.44"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").Value = "'0.0516"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "1.193.24"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'0.0173"
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("D38").Value = "'0.809"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'0.507"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'5.43"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "'0.794"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "1.773.16"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'92.22"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "'54.85"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'7.68"
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D50").Value = "'0.410"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  +0.18%  "
